$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Give "Bibi" her full first name "Bibiana" (last name Blatna unchanged).
$ws.Range("A15").Value = "Bibiana"

# 2. Add a new committee member: Diane Uschner (Roche), with a
#    loc_extended affiliation flag, appended as a new row after row 46.
$ws.Range("A47").Value = "Diane"
$ws.Range("B47").Value = "Uschner"
$ws.Range("C47").Value = "Roche"
$ws.Range("F47").Value = 1

# Leave the selection where it ended up after the last edit.
[void]$ws.Range("G47").Select()
